# Auto-generated script applying numeric corrections to the Leve profit
# calculation columns (H:N) across all eight job sheets, per the scheduled
# market-data refresh run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 3666.6667
$ws.Range("I94").Value = 3600
$ws.Range("K94").Value = 3600
$ws.Range("M94").Value = -3149
$ws.Range("H132").Value = 4352134
$ws.Range("I132").Value = 5004091
$ws.Range("K132").Value = 15012273
$ws.Range("M132").Value = -15009743
$ws.Range("H138").Value = 5004.55
$ws.Range("I138").Value = 3154.5557
$ws.Range("J138").Value = 6518.1816
$ws.Range("K138").Value = 9463.667099999999
$ws.Range("L138").Value = 19554.5448
$ws.Range("M138").Value = -4323.667099999999
$ws.Range("N138").Value = -29834.5448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3667.45
$ws.Range("I32").Value = 2423.2954
$ws.Range("J32").Value = 12791.25
$ws.Range("K32").Value = 2423.2954
$ws.Range("L32").Value = 12791.25
$ws.Range("M32").Value = -2136.2954
$ws.Range("N32").Value = -13365.25
$ws.Range("H74").Value = 1134.6923
$ws.Range("I74").Value = 875.1
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 875.1
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -1.100000000000023
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 1134.6923
$ws.Range("I77").Value = 875.1
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 4375.5
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -7.5
$ws.Range("N77").Value = -18736
$ws.Range("H88").Value = 2000
$ws.Range("I88").Value = 2000
$ws.Range("K88").Value = 2000
$ws.Range("M88").Value = -1594
$ws.Range("H91").Value = 2000
$ws.Range("I91").Value = 2000
$ws.Range("K91").Value = 2000
$ws.Range("M91").Value = -596
$ws.Range("H122").Value = 3856.1538
$ws.Range("I122").Value = 2266.25
$ws.Range("J122").Value = 6400
$ws.Range("K122").Value = 6798.75
$ws.Range("L122").Value = 19200
$ws.Range("M122").Value = -4348.75
$ws.Range("N122").Value = -24100
$ws.Range("H132").Value = 20836954
$ws.Range("I132").Value = 26318956
$ws.Range("K132").Value = 78956868
$ws.Range("M132").Value = -78954338

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1908.1666
$ws.Range("I86").Value = 1403.3334
$ws.Range("J86").Value = 4432.3335
$ws.Range("K86").Value = 1403.3334
$ws.Range("L86").Value = 4432.3335
$ws.Range("M86").Value = -280.3334
$ws.Range("N86").Value = -6678.3335
$ws.Range("H89").Value = 1908.1666
$ws.Range("I89").Value = 1403.3334
$ws.Range("J89").Value = 4432.3335
$ws.Range("K89").Value = 7016.666999999999
$ws.Range("L89").Value = 22161.6675
$ws.Range("M89").Value = -1400.666999999999
$ws.Range("N89").Value = -33393.6675

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2328790.5
$ws.Range("I31").Value = 3127031
$ws.Range("J31").Value = 6636.364
$ws.Range("K31").Value = 3127031
$ws.Range("L31").Value = 6636.364
$ws.Range("M31").Value = -3126736
$ws.Range("N31").Value = -7226.364
$ws.Range("H34").Value = 2328790.5
$ws.Range("I34").Value = 3127031
$ws.Range("J34").Value = 6636.364
$ws.Range("K34").Value = 3127031
$ws.Range("L34").Value = 6636.364
$ws.Range("M34").Value = -3126829
$ws.Range("N34").Value = -7040.364
$ws.Range("H94").Value = 35716620
$ws.Range("I94").Value = 4198
$ws.Range("J94").Value = 38463732
$ws.Range("K94").Value = 4198
$ws.Range("L94").Value = 38463732
$ws.Range("M94").Value = -3747
$ws.Range("N94").Value = -38464634
$ws.Range("H134").Value = 1448.2333
$ws.Range("J134").Value = 2289.4167
$ws.Range("L134").Value = 6868.250100000001
$ws.Range("N134").Value = -11938.2501

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 93.59999999999999
$ws.Range("I2").Value = 70
$ws.Range("J2").Value = 109.333336
$ws.Range("K2").Value = 420
$ws.Range("L2").Value = 656.000016
$ws.Range("M2").Value = -307
$ws.Range("N2").Value = -882.000016
$ws.Range("H34").Value = 17783.334
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 17783.334
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 53350.00199999999
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -53518.00199999999
$ws.Range("H39").Value = 1499.3
$ws.Range("J39").Value = 1750.5
$ws.Range("L39").Value = 5251.5
$ws.Range("N39").Value = -5839.5
$ws.Range("H131").Value = 1296.9788
$ws.Range("I131").Value = 1762
$ws.Range("J131").Value = 1171.2972
$ws.Range("K131").Value = 5286
$ws.Range("L131").Value = 3513.8916
$ws.Range("M131").Value = -246
$ws.Range("N131").Value = -13593.8916
$ws.Range("H136").Value = 2182.1724
$ws.Range("I136").Value = 1510.7
$ws.Range("J136").Value = 3674.3333
$ws.Range("K136").Value = 4532.1
$ws.Range("L136").Value = 11022.9999
$ws.Range("M136").Value = 567.8999999999996
$ws.Range("N136").Value = -21222.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 77.75
$ws.Range("I2").Value = 65.71429000000001
$ws.Range("J2").Value = 87.111115
$ws.Range("K2").Value = 65.71429000000001
$ws.Range("L2").Value = 87.111115
$ws.Range("M2").Value = 47.28570999999999
$ws.Range("N2").Value = -313.111115
$ws.Range("H122").Value = 5355.75
$ws.Range("I122").Value = 4125.875
$ws.Range("J122").Value = 6175.6665
$ws.Range("K122").Value = 12377.625
$ws.Range("L122").Value = 18526.9995
$ws.Range("M122").Value = -9927.625
$ws.Range("N122").Value = -23426.9995
$ws.Range("H132").Value = 2826.0232
$ws.Range("I132").Value = 2271.1482
$ws.Range("K132").Value = 6813.444600000001
$ws.Range("M132").Value = -4283.444600000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 601.6923
$ws.Range("I16").Value = 585.6667
$ws.Range("J16").Value = 637.75
$ws.Range("K16").Value = 585.6667
$ws.Range("L16").Value = 637.75
$ws.Range("M16").Value = -415.6667
$ws.Range("N16").Value = -977.75
$ws.Range("H22").Value = 142859070
$ws.Range("I22").Value = 500000500
$ws.Range("J22").Value = 2496
$ws.Range("K22").Value = 500000500
$ws.Range("L22").Value = 2496
$ws.Range("M22").Value = -500000205
$ws.Range("N22").Value = -3086
$ws.Range("H27").Value = 142859070
$ws.Range("I27").Value = 500000500
$ws.Range("J27").Value = 2496
$ws.Range("K27").Value = 500000500
$ws.Range("L27").Value = 2496
$ws.Range("M27").Value = -500000393
$ws.Range("N27").Value = -2710
$ws.Range("H40").Value = 2595.111
$ws.Range("I40").Value = 1499.5
$ws.Range("K40").Value = 1499.5
$ws.Range("M40").Value = -1363.5
$ws.Range("H100").Value = 2153.3333
$ws.Range("I100").Value = 1380
$ws.Range("J100").Value = 2705.7144
$ws.Range("K100").Value = 1380
$ws.Range("L100").Value = 2705.7144
$ws.Range("M100").Value = -839
$ws.Range("N100").Value = -3787.7144
$ws.Range("H132").Value = 3130.0334
$ws.Range("I132").Value = 1990.7
$ws.Range("J132").Value = 3699.7
$ws.Range("K132").Value = 5972.1
$ws.Range("L132").Value = 11099.1
$ws.Range("M132").Value = -3442.1
$ws.Range("N132").Value = -16159.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1500
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1500
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608
$ws.Range("H122").Value = 264757.7
$ws.Range("I122").Value = 313905.12
$ws.Range("J122").Value = 2638
$ws.Range("K122").Value = 941715.36
$ws.Range("L122").Value = 7914
$ws.Range("M122").Value = -939265.36
$ws.Range("N122").Value = -12814
$ws.Range("H126").Value = 6252320
$ws.Range("I126").Value = 2025.9166
$ws.Range("J126").Value = 25003202
$ws.Range("K126").Value = 6077.7498
$ws.Range("L126").Value = 75009606
$ws.Range("M126").Value = -3607.7498
$ws.Range("N126").Value = -75014546
$ws.Range("H132").Value = 275940.62
$ws.Range("I132").Value = 402648.2
$ws.Range("J132").Value = 11966.5
$ws.Range("K132").Value = 1207944.6
$ws.Range("L132").Value = 35899.5
$ws.Range("M132").Value = -1205414.6
$ws.Range("N132").Value = -40959.5
